$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target table now covers a 4 (sending cluster) x 5 (target cluster) grid of
# Epo -> Epor rows (rows 2-21) because a new target cluster,
# "Inflammatory-Mac", was added to the NATMI run. Every row's metrics were
# recomputed against the new TPM values, and within each sending-cluster
# block the existing target-cluster rows were reordered, so every data cell
# A2:T21 is rewritten explicitly below.
$data = @(
  @{row=2; vals=@("ECs","Epo","Epor","ECs",1,0.3333333333333333,0.07381966666666666,0.221459,0.1284640970637474,0.1284640970637474,3,1,1.374029666666667,4.122089,0.4703393650030796,0.4703393650030795,0.1014304119834444,0.912873707851,0.06042172183865692,0.06042172183865691)}
  @{row=3; vals=@("ECs","Epo","Epor","FAPs",1,0.3333333333333333,0.07381966666666666,0.221459,0.1284640970637474,0.1284640970637474,2,0.6666666666666666,0.3661236666666667,1.098371,0.1253265319302416,0.1253265319302416,0.02702712703211111,0.243244143289,0.01609995976254939,0.01609995976254939)}
  @{row=4; vals=@("ECs","Epo","Epor","Inflammatory-Mac",1,0.3333333333333333,0.07381966666666666,0.221459,0.1284640970637474,0.1284640970637474,2,0.6666666666666666,0.5727720000000001,1.718316,0.1960636115121803,0.1960636115121803,0.042281838116,0.380536543044,0.02518713481996959,0.02518713481996959)}
  @{row=5; vals=@("ECs","Epo","Epor","MuSCs",1,0.3333333333333333,0.07381966666666666,0.221459,0.1284640970637474,0.1284640970637474,3,1,0.4303706666666667,1.291112,0.1473187013254338,0.1473187013254338,0.03176981915644444,0.285928372408,0.01892516394637574,0.01892516394637574)}
  @{row=6; vals=@("FAPs","Epo","Epor","Resolving-Mac",1,0.3333333333333333,0.07381966666666666,0.221459,0.1284640970637474,0.1284640970637474,2,0.6666666666666666,0.178062,0.534186,0.06095179022906471,0.0609517902290647,0.013144477486,0.118300297374,0.007830116696195738,0.007830116696195736)}
  @{row=7; vals=@("FAPs","Epo","Epor","ECs",2,0.6666666666666666,0.1698756666666667,0.5096270000000001,0.2956247991470493,0.2956247991470493,3,1,1.374029666666667,4.122089,0.4703393650030796,0.4703393650030795,0.2334142056447778,2.100727850803,0.1390439803099861,0.1390439803099861)}
  @{row=8; vals=@("FAPs","Epo","Epor","FAPs",2,0.6666666666666666,0.1698756666666667,0.5096270000000001,0.2956247991470493,0.2956247991470493,2,0.6666666666666666,0.3661236666666667,1.098371,0.1253265319302416,0.1253265319302416,0.06219550195744445,0.559759517617,0.03704963082967393,0.03704963082967393)}
  @{row=9; vals=@("FAPs","Epo","Epor","Inflammatory-Mac",2,0.6666666666666666,0.1698756666666667,0.5096270000000001,0.2956247991470493,0.2956247991470493,2,0.6666666666666666,0.5727720000000001,1.718316,0.1960636115121803,0.1960636115121803,0.09730002534800002,0.8757002281320002,0.0579612657733334,0.0579612657733334)}
  @{row=10; vals=@("FAPs","Epo","Epor","MuSCs",2,0.6666666666666666,0.1698756666666667,0.5096270000000001,0.2956247991470493,0.2956247991470493,3,1,0.4303706666666667,1.291112,0.1473187013254338,0.1473187013254338,0.07310950391377778,0.6579855352240001,0.04355106148993552,0.04355106148993552)}
  @{row=11; vals=@("FAPs","Epo","Epor","Resolving-Mac",2,0.6666666666666666,0.1698756666666667,0.5096270000000001,0.2956247991470493,0.2956247991470493,2,0.6666666666666666,0.178062,0.534186,0.06095179022906471,0.0609517902290647,0.03024840095800001,0.2722356086220001,0.01801886074412034,0.01801886074412034)}
  @{row=12; vals=@("MuSCs","Epo","Epor","ECs",2,0.6666666666666666,0.07389766666666667,0.221693,0.1285998359531712,0.1285998359531712,3,1,1.374029666666667,4.122089,0.4703393650030796,0.4703393650030795,0.1015375862974444,0.913838276677,0.06048556518171477,0.06048556518171476)}
  @{row=13; vals=@("MuSCs","Epo","Epor","FAPs",2,0.6666666666666666,0.07389766666666667,0.221693,0.1285998359531712,0.1285998359531712,2,0.6666666666666666,0.3661236666666667,1.098371,0.1253265319302416,0.1253265319302416,0.02705568467811111,0.243501162103,0.01611697144680894,0.01611697144680894)}
  @{row=14; vals=@("MuSCs","Epo","Epor","Inflammatory-Mac",2,0.6666666666666666,0.07389766666666667,0.221693,0.1285998359531712,0.1285998359531712,2,0.6666666666666666,0.5727720000000001,1.718316,0.1960636115121803,0.1960636115121803,0.04232651433200001,0.380938628988,0.02521374827685268,0.02521374827685268)}
  @{row=15; vals=@("MuSCs","Epo","Epor","MuSCs",2,0.6666666666666666,0.07389766666666667,0.221693,0.1285998359531712,0.1285998359531712,3,1,0.4303706666666667,1.291112,0.1473187013254338,0.1473187013254338,0.03180338806844445,0.286230492616,0.01894516082328502,0.01894516082328502)}
  @{row=16; vals=@("MuSCs","Epo","Epor","Resolving-Mac",2,0.6666666666666666,0.07389766666666667,0.221693,0.1285998359531712,0.1285998359531712,2,0.6666666666666666,0.178062,0.534186,0.06095179022906471,0.0609517902290647,0.013158366322,0.118425296898,0.007838390224509826,0.007838390224509826)}
  @{row=17; vals=@("Resolving-Mac","Epo","Epor","ECs",3,1,0.2570396666666667,0.771119,0.447311267836032,0.447311267836032,3,1,1.374029666666667,4.122089,0.4703393650030796,0.4703393650030795,0.3531801275101111,3.178621147591,0.2103880976727218,0.2103880976727217)}
  @{row=18; vals=@("Resolving-Mac","Epo","Epor","FAPs",3,1,0.2570396666666667,0.771119,0.447311267836032,0.447311267836032,2,0.6666666666666666,0.3661236666666667,1.098371,0.1253265319302416,0.1253265319302416,0.09410830523877778,0.846974747149,0.05605996989120931,0.05605996989120931)}
  @{row=19; vals=@("Resolving-Mac","Epo","Epor","Inflammatory-Mac",3,1,0.2570396666666667,0.771119,0.447311267836032,0.447311267836032,2,0.6666666666666666,0.5727720000000001,1.718316,0.1960636115121803,0.1960636115121803,0.147225123956,1.325026115604,0.08770146264202461,0.08770146264202461)}
  @{row=20; vals=@("Resolving-Mac","Epo","Epor","MuSCs",3,1,0.2570396666666667,0.771119,0.447311267836032,0.447311267836032,3,1,0.4303706666666667,1.291112,0.1473187013254338,0.1473187013254338,0.1106223327031111,0.995600994328,0.06589731506583753,0.06589731506583753)}
  @{row=21; vals=@("Resolving-Mac","Epo","Epor","Resolving-Mac",3,1,0.2570396666666667,0.771119,0.447311267836032,0.447311267836032,2,0.6666666666666666,0.178062,0.534186,0.06095179022906471,0.0609517902290647,0.04576899712600001,0.411920974134,0.0272644225642388,0.0272644225642388)}
)

foreach ($item in $data) {
  $r = $item.row
  $vals = $item.vals
  for ($i = 0; $i -lt $vals.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item($r, $col).Value = $vals[$i]
  }
}
